# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (same column layout/styling) to create
#    a new "2022-Q1" sheet, inserted right before the "总计" (totals) sheet.
# 2. Update the fund-scale/position figures on the new sheet to the
#    2022-Q1 numbers (fund code/name/rank stay the same).
# 3. Insert a new leading row on the "总计" sheet summarising 2022-Q1,
#    pushing the existing 2021-Q4 / 2021-Q3 rows down and renumbering the
#    index column.
#
# NOTE: worksheet handles in this host re-resolve by tab position, so any
# variable captured before a sheet-list mutation (Add/Copy/Delete/Move)
# can silently start pointing at a different sheet afterwards. We always
# re-fetch sheets by name immediately before using them post-mutation.

$wb = $excel.ActiveWorkbook

# --- 1) Clone "2021-Q4" -> "2022-Q1", placed just before "总计" ---------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetBefore = $wb.Worksheets.Item("总计")
$q4Sheet.Copy($totalSheetBefore)

# Re-fetch: the sheet list changed, so grab fresh handles by name/value.
$q1Sheet = $wb.Worksheets.Item("2021-Q4 (2)")
$q1Sheet.Name = "2022-Q1"

# --- 2) Overwrite the figures that changed for 2022-Q1 ------------------
# Fund scale / total stock position / position share / held value are
# stored as text in this workbook (matching the other quarter sheets), so
# force text with a leading apostrophe to stop the numeric auto-convert.
$q1Sheet.Range("D2").Value = "'6.05"
$q1Sheet.Range("E2").Value = "'99.49"
$q1Sheet.Range("F2").Value = "'8.08"
$q1Sheet.Range("G2").Value = "'0.4888"
# H2 (仓位排名) keeps its old value of 3, B2/C2 (code/name) are unchanged.

# --- 3) Shift "总计" rows down and insert the 2022-Q1 summary row -------
# Values are written as the known literal targets (rather than copied
# cell-to-cell through Value2) to avoid floating-point round-trip noise.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.17
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.26

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.49

# Restore the originally active tab (sheet copy/add operations move the
# selection onto the new sheet as a side effect).
$wb.Worksheets.Item("2021-Q3").Activate()
